$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string into a cell without Excel silently
# re-interpreting number-looking text (e.g. "580.42", "1.00", "0.0000178")
# as a numeric value. Temporarily force Text format, assign, then restore
# the cells original style so formatting in the file is unaffected.
function Set-LiteralText($rangeAddr, $text) {
    $cell = $ws.Range($rangeAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-LiteralText "D2" "65.197.94"
Set-LiteralText "E2" "  -0.63%  "
Set-LiteralText "D3" "3.403.79"
Set-LiteralText "E3" "  -3.41%  "
Set-LiteralText "E4" "  -0.01%  "
Set-LiteralText "D5" "580.42"
Set-LiteralText "E5" "  -3.26%  "
Set-LiteralText "D6" "135.97"
Set-LiteralText "E6" "  -5.12%  "
Set-LiteralText "E7" "  -0.04%  "
Set-LiteralText "D8" "3.400.85"
Set-LiteralText "E8" "  -3.47%  "
Set-LiteralText "D9" "0.493"
Set-LiteralText "E9" "  -2.08%  "
Set-LiteralText "D10" "7.14"
Set-LiteralText "E10" "  -8.35%  "
Set-LiteralText "D11" "0.120"
Set-LiteralText "E11" "  -10.74%  "
Set-LiteralText "D12" "0.372"
Set-LiteralText "E12" "  -7.71%  "
Set-LiteralText "D13" "3.982.13"
Set-LiteralText "E13" "  -3.52%  "
Set-LiteralText "D14" "0.0000178"
Set-LiteralText "E14" "  -10.20%  "
Set-LiteralText "E15" "  -1.61%  "
Set-LiteralText "B16" "WrappedBTC"
Set-LiteralText "C16" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-LiteralText "D16" "65.148.21"
Set-LiteralText "E16" "  -0.70%  "
Set-LiteralText "B17" "WrappedEther"
Set-LiteralText "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-LiteralText "D17" "3.396.33"
Set-LiteralText "E17" "  -3.71%  "
Set-LiteralText "D18" "25.92"
Set-LiteralText "E18" "  -9.21%  "
Set-LiteralText "D19" "9.53"
Set-LiteralText "E19" "  -12.86%  "
Set-LiteralText "D20" "5.84"
Set-LiteralText "E20" "  -5.75%  "
Set-LiteralText "D21" "13.52"
Set-LiteralText "E21" "  -5.40%  "
Set-LiteralText "D22" "380.09"
Set-LiteralText "E22" "  -8.40%  "
Set-LiteralText "D23" "0.549"
Set-LiteralText "E23" "  -8.22%  "
Set-LiteralText "D24" "1.00"
Set-LiteralText "E24" "  +0.10%  "
Set-LiteralText "D25" "71.86"
Set-LiteralText "E25" "  -7.08%  "
Set-LiteralText "D26" "3.538.00"
Set-LiteralText "E26" "  -3.54%  "
Set-LiteralText "E27" "  -9.55%  "
Set-LiteralText "D28" "0.998"
Set-LiteralText "E28" "  -0.26%  "
Set-LiteralText "D29" "7.05"
Set-LiteralText "E29" "  -8.41%  "
Set-LiteralText "D30" "2.20"
Set-LiteralText "E30" "  -9.78%  "
Set-LiteralText "D31" "8.00"
Set-LiteralText "E31" "  -10.10%  "
Set-LiteralText "D32" "3.413.29"
Set-LiteralText "E32" "  -3.21%  "
Set-LiteralText "E33" "  +0.01%  "
Set-LiteralText "D34" "0.142"
Set-LiteralText "E34" "  -6.82%  "
Set-LiteralText "D35" "22.74"
Set-LiteralText "E35" "  -6.39%  "
Set-LiteralText "D36" "169.55"
Set-LiteralText "E36" "  -2.79%  "
Set-LiteralText "D37" "6.68"
Set-LiteralText "E37" "  -11.31%  "
Set-LiteralText "D38" "1.15"
Set-LiteralText "E38" "  -10.39%  "
Set-LiteralText "D39" "1.46"
Set-LiteralText "E39" "  -7.13%  "
Set-LiteralText "D40" "4.68"
Set-LiteralText "E40" "  -11.31%  "
Set-LiteralText "D41" "0.0752"
Set-LiteralText "E41" "  -8.23%  "
Set-LiteralText "D42" "0.804"
Set-LiteralText "E42" "  -6.02%  "
Set-LiteralText "D43" "43.24"
Set-LiteralText "E43" "  -4.45%  "
Set-LiteralText "D44" "0.999"
Set-LiteralText "E44" "  -0.13%  "
Set-LiteralText "D45" "4.33"
Set-LiteralText "E45" "  -14.62%  "
Set-LiteralText "D46" "1.59"
Set-LiteralText "E46" "  -9.39%  "
Set-LiteralText "D47" "1.10"
Set-LiteralText "E47" "  +1.15%  "
Set-LiteralText "D48" "21.93"
Set-LiteralText "E48" "  -3.23%  "
Set-LiteralText "D49" "6.45"
Set-LiteralText "E49" "  -8.84%  "
Set-LiteralText "D50" "2.05"
Set-LiteralText "E50" "  -13.24%  "
Set-LiteralText "D51" "2.159.87"
Set-LiteralText "E51" "  -7.95%  "
